# Swap the "filter" test step (currently under TC2, row 20) with the
# "cancel" test step (currently under TC3, row 28), so that after the
# edit TC2 holds the cancellation scenario and TC3 holds the filter
# scenario (the TC labels themselves, in B15/B23, stay put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$filterStep   = $ws.Range("B20").Value2
$filterResult = $ws.Range("D20").Value2

$cancelStep   = $ws.Range("B28").Value2
$cancelResult = $ws.Range("D28").Value2

$ws.Range("B20").Value2 = $cancelStep
$ws.Range("D20").Value2 = $cancelResult

$ws.Range("B28").Value2 = $filterStep
$ws.Range("D28").Value2 = $filterResult
